$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update row 2 values to reflect the latest stock ticker data.
# Briefly force text formatting so values like "$68.93" are stored as
# literal text rather than being reinterpreted as numeric/currency
# values, then restore the original (unstyled) cell style.
$ws.Range("B2").NumberFormat = "@"
$ws.Range("B2").Value = "$68.93"
$ws.Range("B2").Style = "Normal"

$ws.Range("C2").NumberFormat = "@"
$ws.Range("C2").Value = "-0.02(0.03%) 1D"
$ws.Range("C2").Style = "Normal"

# Remove rows 3-8 (previously hard-coded stock list), leaving only the
# single ticker-driven row
$ws.Range("A3:D8").EntireRow.Delete()
